# refactor currency conversion, now explicit source and target amounts
$wb = $excel.ActiveWorkbook

# --- currency_conversions sheet: rewrite header row with explicit
#     source_amount / source_fees / source_currency / target_amount /
#     target_fees / target_currency / comment columns (drops the old
#     single "foreign_amount" column) ---
$wsCur = $wb.Worksheets.Item(5)
$wsCur.Range("A1").Value = "date"
$wsCur.Range("B1").Value = "source_amount"
$wsCur.Range("C1").Value = "source_fees"
$wsCur.Range("D1").Value = "source_currency"
$wsCur.Range("E1").Value = "target_amount"
$wsCur.Range("F1").Value = "target_fees"
$wsCur.Range("G1").Value = "target_currency"
$wsCur.Range("H1").Value = "comment"
$wsCur.Range("A1").Copy()
$wsCur.Range("B1:H1").PasteSpecial(-4122)

# --- buy_orders: drop per-row formatting on the data rows and remove the
#     trailing empty column G cells ---
$wsBuy = $wb.Worksheets.Item(1)
$wsBuy.Range("B2:F3").ClearFormats()
$wsBuy.Range("G2:G3").Clear()

# --- sell_orders: same cleanup as buy_orders (single data row) ---
$wsSell = $wb.Worksheets.Item(4)
$wsSell.Range("B2:F2").ClearFormats()
$wsSell.Range("G2").Clear()

# --- money_transfers: remove the stray formatted row 20 ---
$wsMoney = $wb.Worksheets.Item(3)
$wsMoney.Rows.Item(20).Delete()

# --- make currency_conversions the active tab/sheet, as in the source file ---
$wsCur.Activate()
